$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 172.5
$ws.Range("I19").Value = 288.33334
$ws.Range("J19").Value = 56.666668
$ws.Range("K19").Value = 288.33334
$ws.Range("L19").Value = 56.666668
$ws.Range("M19").Value = -113.33334
$ws.Range("N19").Value = -406.666668
$ws.Range("H40").Value = 2924.625
$ws.Range("J40").Value = 2927.9
$ws.Range("L40").Value = 2927.9
$ws.Range("N40").Value = -3277.9
$ws.Range("H70").Value = 2170.125
$ws.Range("J70").Value = 2999.75
$ws.Range("L70").Value = 8999.25
$ws.Range("N70").Value = -9539.25
$ws.Range("H73").Value = 2170.125
$ws.Range("J73").Value = 2999.75
$ws.Range("L73").Value = 8999.25
$ws.Range("N73").Value = -10871.25
$ws.Range("H80").Value = 1816.6154
$ws.Range("I80").Value = 1731.6364
$ws.Range("K80").Value = 5194.9092
$ws.Range("M80").Value = -4196.9092
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H83").Value = 1816.6154
$ws.Range("I83").Value = 1731.6364
$ws.Range("K83").Value = 15584.7276
$ws.Range("M83").Value = -10592.7276
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H86").Value = 71431610
$ws.Range("I86").Value = 111114190
$ws.Range("J86").Value = 2960.8
$ws.Range("K86").Value = 111114190
$ws.Range("L86").Value = 2960.8
$ws.Range("M86").Value = -111113067
$ws.Range("N86").Value = -5206.8
$ws.Range("H89").Value = 71431610
$ws.Range("I89").Value = 111114190
$ws.Range("J89").Value = 2960.8
$ws.Range("K89").Value = 555570950
$ws.Range("L89").Value = 14804
$ws.Range("M89").Value = -555565334
$ws.Range("N89").Value = -26036
$ws.Range("H125").Value = 3366.8
$ws.Range("I125").Value = 1500
$ws.Range("J125").Value = 3574.2222
$ws.Range("K125").Value = 13500
$ws.Range("L125").Value = 32167.9998
$ws.Range("M125").Value = -11040
$ws.Range("N125").Value = -37087.99980000001
$ws.Range("H131").Value = 7609
$ws.Range("I131").Value = 6846.1816
$ws.Range("K131").Value = 20538.5448
$ws.Range("M131").Value = -15498.5448
$ws.Range("H136").Value = 143747
$ws.Range("J136").Value = 143747
$ws.Range("L136").Value = 143747
$ws.Range("N136").Value = -153947
$ws.Range("H138").Value = 5431.4033
$ws.Range("J138").Value = 2739.325
$ws.Range("L138").Value = 8217.974999999999
$ws.Range("N138").Value = -18497.975
$ws.Range("H141").Value = 2982.3333
$ws.Range("J141").Value = 2900
$ws.Range("L141").Value = 8700
$ws.Range("N141").Value = -19060

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 333333700
$ws.Range("I30").Value = 500000200
$ws.Range("J30").Value = 650
$ws.Range("K30").Value = 500000200
$ws.Range("L30").Value = 650
$ws.Range("M30").Value = -500000050
$ws.Range("N30").Value = -950
$ws.Range("H32").Value = 182954.75
$ws.Range("I32").Value = 198069.75
$ws.Range("J32").Value = 28781.8
$ws.Range("K32").Value = 198069.75
$ws.Range("L32").Value = 28781.8
$ws.Range("M32").Value = -197782.75
$ws.Range("N32").Value = -29355.8
$ws.Range("H45").Value = 73788.14
$ws.Range("I45").Value = 79348.84
$ws.Range("K45").Value = 79348.84
$ws.Range("M45").Value = -78971.84
$ws.Range("H110").Value = 653.8570999999999
$ws.Range("I110").Value = 653.8570999999999
$ws.Range("K110").Value = 653.8570999999999
$ws.Range("M110").Value = 1391.1429
$ws.Range("H124").Value = 20714.5
$ws.Range("J124").Value = 20714.5
$ws.Range("L124").Value = 20714.5
$ws.Range("N124").Value = -30534.5
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4778.909
$ws.Range("I86").Value = 1830.5834
$ws.Range("K86").Value = 1830.5834
$ws.Range("M86").Value = -707.5834
$ws.Range("H89").Value = 4778.909
$ws.Range("I89").Value = 1830.5834
$ws.Range("K89").Value = 9152.916999999999
$ws.Range("M89").Value = -3536.916999999999
$ws.Range("H105").Value = 8336.857
$ws.Range("I105").Value = 10103.692
$ws.Range("J105").Value = 5465.75
$ws.Range("K105").Value = 10103.692
$ws.Range("L105").Value = 5465.75
$ws.Range("M105").Value = -8356.691999999999
$ws.Range("N105").Value = -8959.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 11150
$ws.Range("I6").Value = 17500
$ws.Range("J6").Value = 4800
$ws.Range("K6").Value = 17500
$ws.Range("L6").Value = 4800
$ws.Range("M6").Value = -17387
$ws.Range("N6").Value = -5026
$ws.Range("H58").Value = 2717
$ws.Range("I58").Value = 1674.125
$ws.Range("J58").Value = 3551.3
$ws.Range("K58").Value = 1674.125
$ws.Range("L58").Value = 3551.3
$ws.Range("M58").Value = -1471.125
$ws.Range("N58").Value = -3957.3
$ws.Range("H107").Value = 1407.5758
$ws.Range("J107").Value = 1274.7059
$ws.Range("L107").Value = 1274.7059
$ws.Range("N107").Value = -5114.7059
$ws.Range("H122").Value = 3410.8333
$ws.Range("I122").Value = 2860.1333
$ws.Range("K122").Value = 8580.3999
$ws.Range("M122").Value = -6130.3999
$ws.Range("H132").Value = 16670430
$ws.Range("J132").Value = 37040350
$ws.Range("L132").Value = 111121050
$ws.Range("N132").Value = -111126110
$ws.Range("H136").Value = 2717
$ws.Range("I136").Value = 1674.125
$ws.Range("J136").Value = 3551.3
$ws.Range("K136").Value = 5022.375
$ws.Range("L136").Value = 10653.9
$ws.Range("M136").Value = -2472.375
$ws.Range("N136").Value = -15753.9
$ws.Range("H138").Value = 95249.5
$ws.Range("I138").Value = 95500
$ws.Range("J138").Value = 94999
$ws.Range("K138").Value = 95500
$ws.Range("L138").Value = 94999
$ws.Range("M138").Value = -90360
$ws.Range("N138").Value = -105279

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 564.75
$ws.Range("I5").Value = 557.9
$ws.Range("K5").Value = 1673.7
$ws.Range("M5").Value = -1561.7
$ws.Range("H107").Value = 37037396
$ws.Range("I107").Value = 758.6
$ws.Range("K107").Value = 2275.8
$ws.Range("M107").Value = -355.8000000000002
$ws.Range("H135").Value = 564.75
$ws.Range("I135").Value = 557.9
$ws.Range("K135").Value = 5021.099999999999
$ws.Range("M135").Value = -2486.099999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10835.565
$ws.Range("I70").Value = 12043.667
$ws.Range("K70").Value = 12043.667
$ws.Range("M70").Value = -11773.667
$ws.Range("H73").Value = 10835.565
$ws.Range("I73").Value = 12043.667
$ws.Range("K73").Value = 12043.667
$ws.Range("M73").Value = -11107.667
$ws.Range("H80").Value = 146214.33
$ws.Range("I80").Value = 144357.2
$ws.Range("J80").Value = 155500
$ws.Range("K80").Value = 144357.2
$ws.Range("L80").Value = 155500
$ws.Range("M80").Value = -143359.2
$ws.Range("N80").Value = -157496
$ws.Range("H83").Value = 146214.33
$ws.Range("I83").Value = 144357.2
$ws.Range("J83").Value = 155500
$ws.Range("K83").Value = 721786
$ws.Range("L83").Value = 777500
$ws.Range("M83").Value = -716794
$ws.Range("N83").Value = -787484
$ws.Range("H122").Value = 2783.1667
$ws.Range("I122").Value = 1362
$ws.Range("K122").Value = 4086
$ws.Range("M122").Value = -1636
$ws.Range("H126").Value = 1980.5714
$ws.Range("I126").Value = 1862.6154
$ws.Range("K126").Value = 5587.8462
$ws.Range("M126").Value = -3117.8462
$ws.Range("H132").Value = 14710313
$ws.Range("I132").Value = 3856
$ws.Range("K132").Value = 11568
$ws.Range("M132").Value = -9038
$ws.Range("H136").Value = 30324.445
$ws.Range("J136").Value = 30324.445
$ws.Range("L136").Value = 90973.33499999999
$ws.Range("N136").Value = -96073.33499999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1868913
$ws.Range("I40").Value = 2241576
$ws.Range("K40").Value = 2241576
$ws.Range("M40").Value = -2241440
$ws.Range("H46").Value = 5144
$ws.Range("J46").Value = 1886.75
$ws.Range("L46").Value = 1886.75
$ws.Range("N46").Value = -2262.75
$ws.Range("H94").Value = 59523.43
$ws.Range("H122").Value = 5106.4644
$ws.Range("I122").Value = 3499.625
$ws.Range("J122").Value = 5749.2
$ws.Range("K122").Value = 10498.875
$ws.Range("L122").Value = 17247.6
$ws.Range("M122").Value = -8048.875
$ws.Range("N122").Value = -22147.6
$ws.Range("H132").Value = 13113.077
$ws.Range("I132").Value = 4500
$ws.Range("K132").Value = 13500
$ws.Range("M132").Value = -10970

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 49999.25
$ws.Range("I14").Value = 49999
$ws.Range("K14").Value = 49999
$ws.Range("M14").Value = -49831
$ws.Range("H122").Value = 21273.312
$ws.Range("I122").Value = 22591.533
$ws.Range("K122").Value = 67774.599
$ws.Range("M122").Value = -65324.599
$ws.Range("H126").Value = 3063.16
$ws.Range("I126").Value = 2906.889
$ws.Range("J126").Value = 3465
$ws.Range("K126").Value = 8720.667000000001
$ws.Range("L126").Value = 10395
$ws.Range("M126").Value = -6250.667000000001
$ws.Range("N126").Value = -15335
$ws.Range("H132").Value = 35002.266
$ws.Range("I132").Value = 51391.2
$ws.Range("J132").Value = 2224.4
$ws.Range("K132").Value = 154173.6
$ws.Range("L132").Value = 6673.200000000001
$ws.Range("M132").Value = -151643.6
$ws.Range("N132").Value = -11733.2
